$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the per-dataset error-summary formulas (columns Q:Y, rows 2-13) and
# replace them with a placeholder single-space text in column Q, leaving
# R:Y empty -- the detailed error write-up will be filled in manually later.
for ($r = 2; $r -le 13; $r++) {
    $ws.Range("Q$r").Value = " "
    $ws.Range("R$r`:Y$r").ClearContents()
}

# Re-merge the "Average" row labels so the merge-cell bookkeeping matches
# the state Excel produced after the above edits.
foreach ($ref in @("A5:D5","A10:D10","A16:D16","A21:D21","A58:D58")) {
    $ws.Range($ref).UnMerge()
    $ws.Range($ref).Merge()
}

# Update the saved selection/view state to match the author's last position.
$ws.Range("U17").Select()
